# Updated cryptos list on Fri Mar 31 13:31:18 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be stored as text so numeric-looking strings
    # (e.g. "0.9987", "41.69") are not auto-converted to numbers and
    # keep their exact original formatting (trailing zeros, etc.).
    $c = $ws.Range($range)
    $c.NumberFormat = "@"
    $c.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "28.235.35"
$ws.Range("E2").Value = "  -1.37%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.815.83"
$ws.Range("E3").Value = "  +0.54%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.9987"
$ws.Range("E4").Value = "  -0.26%  "

# Row 5 - BNB
Set-TextValue "D5" "316.81"
$ws.Range("E5").Value = "  -0.45%  "

# Row 6 - USDC
Set-TextValue "D6" "0.9982"
$ws.Range("E6").Value = "  -0.34%  "

# Row 7 - XRP
Set-TextValue "D7" "0.5317"
$ws.Range("E7").Value = "  -2.55%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.3976"
$ws.Range("E8").Value = "  +4.58%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.07538"
$ws.Range("E9").Value = "  +0.15%  "

# Row 10 - OKB
Set-TextValue "D10" "41.69"
$ws.Range("E10").Value = "  -1.57%  "

# Row 11 - Polygon (only Volume changes)
$ws.Range("E11").Value = "  -1.73%  "

# Row 12 - swaps from BinanceUSD to Chainlink
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D12" "7.614"
$ws.Range("E12").Value = "  +3.44%  "

# Row 13 - swaps from Chainlink to BinanceUSD
$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D13" "0.9973"
$ws.Range("E13").Value = "  -0.48%  "

# Row 14 - Polkadot
Set-TextValue "D14" "6.243"
$ws.Range("E14").Value = "  +1.23%  "

# Row 15 - Solana
Set-TextValue "D15" "20.53"
$ws.Range("E15").Value = "  -0.77%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "1.804.24"
$ws.Range("E16").Value = "  +0.35%  "

# Row 17 - Litecoin
Set-TextValue "D17" "89.19"
$ws.Range("E17").Value = "  -1.23%  "

# Row 18 - ShibaInu
Set-TextValue "D18" "0.00001065"
$ws.Range("E18").Value = "  -0.23%  "

# Row 19 - TRON (only Volume changes)
$ws.Range("E19").Value = "  +0.50%  "

# Row 20 - Dai
Set-TextValue "D20" "0.9980"
$ws.Range("E20").Value = "  -0.28%  "

# Row 21 - Avalanche
Set-TextValue "D21" "17.34"
$ws.Range("E21").Value = "  -0.85%  "

# Row 22 - Uniswap
Set-TextValue "D22" "6.015"
$ws.Range("E22").Value = "  +0.94%  "

# Row 23 - WrappedBTC
Set-TextValue "D23" "28.269.68"
$ws.Range("E23").Value = "  -1.26%  "

# Row 24 - Cosmos
Set-TextValue "D24" "11.14"
$ws.Range("E24").Value = "  -0.18%  "

# Row 25 - Toncoin
Set-TextValue "D25" "2.090"
$ws.Range("E25").Value = "  +0.53%  "

# Row 26 - Monero
Set-TextValue "D26" "156.12"
$ws.Range("E26").Value = "  -3.23%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "20.36"
$ws.Range("E27").Value = "  -0.87%  "

# Row 28 - WrappedliquidstakedEther2.0
Set-TextValue "D28" "2.015.05"
$ws.Range("E28").Value = "  +0.44%  "

# Row 29 - LidoDAOToken
Set-TextValue "D29" "2.376"
$ws.Range("E29").Value = "  +1.26%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "123.08"
$ws.Range("E30").Value = "  -0.14%  "

# Row 31 - Stellar
Set-TextValue "D31" "0.1095"
$ws.Range("E31").Value = "  +3.13%  "

# Row 32 - ImmutableX (only Volume changes)
$ws.Range("E32").Value = "  -2.68%  "

# Row 33 - HuobiToken
Set-TextValue "D33" "3.669"
$ws.Range("E33").Value = "  -0.31%  "

# Row 34 - swaps from Filecoin to Hedera
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D34" "0.07393"
$ws.Range("E34").Value = "  +11.92%  "

# Row 35 - swaps from Hedera to Filecoin
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D35" "5.560"
$ws.Range("E35").Value = "  -1.71%  "

# Row 36 - Algorand
Set-TextValue "D36" "0.2228"
$ws.Range("E36").Value = "  -1.59%  "

# Row 37 - InternetComputer(DFINITY)
Set-TextValue "D37" "5.175"
$ws.Range("E37").Value = "  +2.30%  "

# Row 38 - VeChain
Set-TextValue "D38" "0.02298"
$ws.Range("E38").Value = "  -0.63%  "

# Row 39 - FraxShare
Set-TextValue "D39" "8.617"
$ws.Range("E39").Value = "  -0.20%  "

# Row 40 - Aptos
Set-TextValue "D40" "11.30"
$ws.Range("E40").Value = "  +0.23%  "

# Row 41 - TheSandbox
Set-TextValue "D41" "0.6200"
$ws.Range("E41").Value = "  -0.63%  "

# Row 42 - TrustWalletToken
Set-TextValue "D42" "1.192"
$ws.Range("E42").Value = "  -0.44%  "

# Row 43 - WEMIXTOKEN
Set-TextValue "D43" "1.403"
$ws.Range("E43").Value = "  -3.31%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "13.38"
$ws.Range("E44").Value = "  -0.11%  "

# Row 45 - PancakeSwap
Set-TextValue "D45" "3.685"
$ws.Range("E45").Value = "  -0.36%  "

# Row 46 - Decentraland
Set-TextValue "D46" "0.5762"
$ws.Range("E46").Value = "  -1.68%  "

# Row 47 - Quant
Set-TextValue "D47" "125.09"
$ws.Range("E47").Value = "  -2.08%  "

# Row 48 - NEARProtocol
Set-TextValue "D48" "1.945"
$ws.Range("E48").Value = "  -1.05%  "

# Row 49 - EOS (only Volume changes)
$ws.Range("E49").Value = "  -0.11%  "

# Row 50 - Cronos
Set-TextValue "D50" "0.06832"
$ws.Range("E50").Value = "  -1.11%  "

# Row 51 - Aave
Set-TextValue "D51" "70.92"
$ws.Range("E51").Value = "  -2.33%  "
